# Add a new "scraped_at" column (D) to the organisatietypes sheet and
# stamp every existing data row with the scrape date "2025-05-27".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Figure out how many data rows currently exist (column C has the last
# populated column before this edit, e.g. A1:C78).
$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$lastDataCol = $used.Columns.Count
$newCol = $lastDataCol + 1

# --- Header cell (row 1): reuse the same bold/bordered/centered style as
# the other header cells (A1:C1) by copying C1's formatting onto the new
# header cell, then set its text.
$headerCell = $ws.Cells.Item(1, $lastDataCol)
$newHeaderCell = $ws.Cells.Item(1, $newCol)
$headerCell.Copy() | Out-Null
$newHeaderCell.PasteSpecial(-4122) | Out-Null
$newHeaderCell.Value = "scraped_at"

# --- Data cells (rows 2..lastRow): write the scrape date as literal text
# (not an Excel-auto-converted date serial number).
$dataRowCount = $lastRow - 1
$dateRange = $ws.Range($ws.Cells.Item(2, $newCol), $ws.Cells.Item($lastRow, $newCol))

$values = New-Object 'object[,]' $dataRowCount,1
for ($i = 0; $i -lt $dataRowCount; $i++) {
    $values[$i, 0] = "2025-05-27"
}

# Format as Text first so the ISO-looking date string is kept verbatim
# instead of being reinterpreted as a date value...
$dateRange.NumberFormat = "@"
$dateRange.Value = $values
# ...then drop back to the default "Normal" style (matches the other,
# unstyled data cells) now that the value is safely stored as text.
$dateRange.Style = "Normal"

"Updated range: " + $dateRange.Address()
"New dimension rows=" + $lastRow + " cols=" + $newCol
